$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 10.47446473282706
$ws.Range("B2").Value = 0.0000000000000001110223024625157
$ws.Range("C2").Value = 0.006439227036331363
$ws.Range("D2").Value = 1.130144420033082
$ws.Range("E2").Value = 1.27722641013191

$ws.Range("A3").Value = 8.270020286634312
$ws.Range("B3").Value = 0.0000000000000001110223024625157
$ws.Range("C3").Value = 0.005957891004521166
$ws.Range("D3").Value = 1.045665455796858
$ws.Range("E3").Value = 1.093416245446851

$ws.Range("A4").Value = 9.807276981582211
$ws.Range("B4").Value = 0.0000000000000001110223024625157
$ws.Range("C4").Value = 0.006272018461888336
$ws.Range("D4").Value = 1.100797755235852
$ws.Range("E4").Value = 1.211755697932291

$ws.Range("A5").Value = 9.364191565654082
$ws.Range("B5").Value = 0.0000000000000001110223024625157
$ws.Range("C5").Value = 0.006594109795157517
$ws.Range("D5").Value = 1.15732778919511
$ws.Range("E5").Value = 1.33940761164324

$ws.Range("A6").Value = 7.668068050422462
$ws.Range("B6").Value = 0.0000000000000001110223024625157
$ws.Range("C6").Value = 0.005562046108038715
$ws.Range("D6").Value = 0.9761909834053566
$ws.Range("E6").Value = 0.9529488360819172

